# Adds a new "Development WebUI" entry as row 17 of the Stundenliste and
# shifts the running total (shared formula already spans C7:C35, so the
# new B17 value automatically propagates through the rest of the column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row: date, hours, running-total formula, activity text.
$ws.Range("A17").Value = 44218
$ws.Range("B17").Value = 3
$ws.Range("C17").Formula = "=C16+B17"
$ws.Range("D17").Value = "Development WebUI"

# Give the new date cell the same number format (date) as the cell above it.
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false | Out-Null

# Restore the selection to where the author last left off editing.
$ws.Range("J14").Select() | Out-Null
